$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(1, 1).Value = "user1"
$ws.Cells.Item(1, 2).Value = "pass1"
$ws.Cells.Item(2, 1).Value = "user2"
$ws.Cells.Item(3, 1).Value = "user3"
$ws.Cells.Item(4, 1).Value = "user4"
$ws.Cells.Item(5, 1).Value = "user5"
$ws.Cells.Item(2, 2).Value = "pass2"
$ws.Cells.Item(3, 2).Value = "pass3"
$ws.Cells.Item(4, 2).Value = "pass4"
$ws.Cells.Item(5, 2).Value = "pass5"

$ws.Range("B1:B5").Select()
